$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-16 Monday" "2024-12-17 Tuesday"

Replace-Text "71×43=3053" "12×14=168"
Replace-Text "31×48=1488" "81×14=1134"
Replace-Text "30×33=990" "89×79=7031"
Replace-Text "89×59=5251" "20×99=1980"
Replace-Text "82×55=4510" "84×74=6216"
Replace-Text "91×58=5278" "67×29=1943"
Replace-Text "70×32=2240" "63×17=1071"
Replace-Text "43×13=559" "25×61=1525"
Replace-Text "44×99=4356" "11×79=869"
Replace-Text "71×25=1775" "44×48=2112"
Replace-Text "90×71=6390" "89×37=3293"
Replace-Text "27×46=1242" "40×26=1040"
Replace-Text "31×46=1426" "20×57=1140"
Replace-Text "65×69=4485" "42×42=1764"
Replace-Text "79×64=5056" "42×36=1512"
Replace-Text "53×62=3286" "74×24=1776"
Replace-Text "90×26=2340" "13×97=1261"
Replace-Text "15×18=270" "47×95=4465"
Replace-Text "13×23=299" "74×90=6660"
Replace-Text "73×88=6424" "20×55=1100"
Replace-Text "65×72=4680" "87×35=3045"
Replace-Text "12×17=204" "31×80=2480"
Replace-Text "59×40=2360" "16×85=1360"
Replace-Text "56×28=1568" "93×12=1116"
Replace-Text "52×58=3016" "76×75=5700"

Write-Output "Done"
